# Atualização de bases das ligas, do dia: 29-02-2024 às 07:50
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Row($r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

function Copy-Row($src, $dst) {
    foreach ($col in $cols) {
        $srcAddr = "$col$src"
        $dstAddr = "$col$dst"
        $v = $ws.Range($srcAddr).Value2
        $ws.Range($dstAddr).Value = $v
    }
}

# Fixture results that were re-matched to the correct id / odds record
# (full swap of columns B..AC between the two rows, A and date stay put)
Swap-Row 72 73
Swap-Row 94 95
Swap-Row 98 99
Swap-Row 128 129
Swap-Row 132 133
Swap-Row 148 149
Swap-Row 175 176
Swap-Row 188 189
Swap-Row 193 194
Swap-Row 222 223
Swap-Row 237 238

# Matches that have now been played - final score & derived odds columns filled in
$ws.Range("H248").Value = 4
$ws.Range("I248").Value = 1
$ws.Range("J248").Value = "H"
$ws.Range("N248").Value = 2.625
$ws.Range("O248").Value = 3.4
$ws.Range("P248").Value = 2.625
$ws.Range("W248").Value = 1.625
$ws.Range("X248").Value = -1
$ws.Range("Y248").Value = -1
$ws.Range("Z248").Value = 0.925
$ws.Range("AA248").Value = -1
$ws.Range("AB248").Value = 0.9750000000000001
$ws.Range("AC248").Value = -1

$ws.Range("H249").Value = 1
$ws.Range("I249").Value = 0
$ws.Range("J249").Value = "H"
$ws.Range("R249").Value = 1.875
$ws.Range("S249").Value = 1.975
$ws.Range("W249").Value = 0.3999999999999999
$ws.Range("X249").Value = -1
$ws.Range("Y249").Value = -1
$ws.Range("Z249").Value = -0.5
$ws.Range("AA249").Value = 0.4875
$ws.Range("AB249").Value = -1
$ws.Range("AC249").Value = 0.925

# New day's fixtures pushed the schedule down one row; last (oldest-added) row drops off
for ($r = 250; $r -le 258; $r++) {
    Copy-Row ($r + 1) $r
}

$ws.Rows(259).Delete()
